# edit.ps1
# Applies the Diebold-Mariano correction described by the commit:
# "Correcion a Diebold Mariano y revision de Cap1"
#
# Updates the four worksheets of the workbook:
#  - Matriz_Resultados : corrected win/lose/tie indicator matrix
#  - P_valores         : corrected DM-test p-values
#  - Estadisticos_DM   : corrected DM-test statistics
#  - Resumen           : corrected summary table (wins/losses/ties/rate/ECRPS)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Matriz_Resultados")
$ws2 = $wb.Worksheets.Item("P_valores")
$ws3 = $wb.Worksheets.Item("Estadisticos_DM")
$ws4 = $wb.Worksheets.Item("Resumen")

# --- Matriz_Resultados ---
$ws1.Range("F2").Value2 = 0
$ws1.Range("G3").Value2 = 0
$ws1.Range("H3").Value2 = 0
$ws1.Range("G4").Value2 = 0
$ws1.Range("H4").Value2 = 0
$ws1.Range("B6").Value2 = 0
$ws1.Range("J6").Value2 = 0
$ws1.Range("C7").Value2 = 0
$ws1.Range("D7").Value2 = 0
$ws1.Range("C8").Value2 = 0
$ws1.Range("D8").Value2 = 0
$ws1.Range("F10").Value2 = 0

# --- P_valores ---
$ws2.Range("C2").Value2 = [double]"1.597491117588667E-05"
$ws2.Range("D2").Value2 = [double]"2.191557348751871E-05"
$ws2.Range("E2").Value2 = 0.0002314696878462641
$ws2.Range("F2").Value2 = 0.004910650961950314
$ws2.Range("G2").Value2 = [double]"1.691900325528906E-05"
$ws2.Range("H2").Value2 = [double]"1.332306989398901E-05"
$ws2.Range("I2").Value2 = [double]"7.223860896177214E-05"
$ws2.Range("J2").Value2 = 0.001286727177699953
$ws2.Range("B3").Value2 = [double]"1.597491117588667E-05"
$ws2.Range("D3").Value2 = [double]"7.695170165522569E-07"
$ws2.Range("E3").Value2 = 0.0001739741554593888
$ws2.Range("F3").Value2 = [double]"2.464474382102289E-05"
$ws2.Range("G3").Value2 = 0.002062180621301035
$ws2.Range("H3").Value2 = 0.001711905625133125
$ws2.Range("I3").Value2 = 0.06581052384841657
$ws2.Range("J3").Value2 = [double]"4.805882292124863E-07"
$ws2.Range("B4").Value2 = [double]"2.191557348751871E-05"
$ws2.Range("C4").Value2 = [double]"7.695170165522569E-07"
$ws2.Range("E4").Value2 = 0.0006159196826533986
$ws2.Range("F4").Value2 = [double]"3.545133697335778E-05"
$ws2.Range("G4").Value2 = 0.009411472722824943
$ws2.Range("H4").Value2 = 0.008119211299410178
$ws2.Range("I4").Value2 = 0.290591759030252
$ws2.Range("J4").Value2 = [double]"6.975348503246437E-07"
$ws2.Range("B5").Value2 = 0.0002314696878462641
$ws2.Range("C5").Value2 = 0.0001739741554593888
$ws2.Range("D5").Value2 = 0.0006159196826533986
$ws2.Range("F5").Value2 = 0.0003787231344722919
$ws2.Range("G5").Value2 = 0.527442073010155
$ws2.Range("H5").Value2 = 0.3544697479583658
$ws2.Range("I5").Value2 = 0.1402225749112975
$ws2.Range("J5").Value2 = 0.0003360638276204142
$ws2.Range("B6").Value2 = 0.004910650961950314
$ws2.Range("C6").Value2 = [double]"2.464474382102289E-05"
$ws2.Range("D6").Value2 = [double]"3.545133697335778E-05"
$ws2.Range("E6").Value2 = 0.0003787231344722919
$ws2.Range("G6").Value2 = [double]"5.998858278100094E-05"
$ws2.Range("H6").Value2 = [double]"4.009675443206184E-05"
$ws2.Range("I6").Value2 = 0.0001074156789391267
$ws2.Range("J6").Value2 = 0.005415387987691478
$ws2.Range("B7").Value2 = [double]"1.691900325528906E-05"
$ws2.Range("C7").Value2 = 0.002062180621301035
$ws2.Range("D7").Value2 = 0.009411472722824943
$ws2.Range("E7").Value2 = 0.527442073010155
$ws2.Range("F7").Value2 = [double]"5.998858278100094E-05"
$ws2.Range("H7").Value2 = 0.406849896997644
$ws2.Range("I7").Value2 = 0.2814154251293135
$ws2.Range("J7").Value2 = [double]"6.305509063775361E-07"
$ws2.Range("B8").Value2 = [double]"1.332306989398901E-05"
$ws2.Range("C8").Value2 = 0.001711905625133125
$ws2.Range("D8").Value2 = 0.008119211299410178
$ws2.Range("E8").Value2 = 0.3544697479583658
$ws2.Range("F8").Value2 = [double]"4.009675443206184E-05"
$ws2.Range("G8").Value2 = 0.406849896997644
$ws2.Range("I8").Value2 = 0.3541757020005531
$ws2.Range("J8").Value2 = [double]"1.029608176494889E-09"
$ws2.Range("B9").Value2 = [double]"7.223860896177214E-05"
$ws2.Range("C9").Value2 = 0.06581052384841657
$ws2.Range("D9").Value2 = 0.290591759030252
$ws2.Range("E9").Value2 = 0.1402225749112975
$ws2.Range("F9").Value2 = 0.0001074156789391267
$ws2.Range("G9").Value2 = 0.2814154251293135
$ws2.Range("H9").Value2 = 0.3541757020005531
$ws2.Range("J9").Value2 = 0.001272840497232419
$ws2.Range("B10").Value2 = 0.001286727177699953
$ws2.Range("C10").Value2 = [double]"4.805882292124863E-07"
$ws2.Range("D10").Value2 = [double]"6.975348503246437E-07"
$ws2.Range("E10").Value2 = 0.0003360638276204142
$ws2.Range("F10").Value2 = 0.005415387987691478
$ws2.Range("G10").Value2 = [double]"6.305509063775361E-07"
$ws2.Range("H10").Value2 = [double]"1.029608176494889E-09"
$ws2.Range("I10").Value2 = 0.001272840497232419

# --- Estadisticos_DM ---
$ws3.Range("C2").Value2 = 6.421081901608537
$ws3.Range("D2").Value2 = 6.232214611954224
$ws3.Range("E2").Value2 = 4.90650264817396
$ws3.Range("F2").Value2 = 3.334758169096592
$ws3.Range("G2").Value2 = 6.386564016778547
$ws3.Range("H2").Value2 = 6.530861996158466
$ws3.Range("I2").Value2 = 5.54463229587341
$ws3.Range("J2").Value2 = 4.011492991243025
$ws3.Range("B3").Value2 = -6.421081901608537
$ws3.Range("D3").Value2 = -8.402791546320882
$ws3.Range("E3").Value2 = -5.060307348339624
$ws3.Range("F3").Value2 = -6.162846263287839
$ws3.Range("G3").Value2 = -3.771893612103899
$ws3.Range("H3").Value2 = -3.866228331476032
$ws3.Range("I3").Value2 = -1.995625243972477
$ws3.Range("J3").Value2 = -8.741976627948864
$ws3.Range("B4").Value2 = -6.232214611954224
$ws3.Range("C4").Value2 = 8.402791546320882
$ws3.Range("E4").Value2 = -4.390548150434426
$ws3.Range("F4").Value2 = -5.950397803041026
$ws3.Range("G4").Value2 = -3.0074511223745
$ws3.Range("H4").Value2 = -3.081890267858277
$ws3.Range("I4").Value2 = -1.098328329486271
$ws3.Range("J4").Value2 = -8.472778910051924
$ws3.Range("B5").Value2 = -4.90650264817396
$ws3.Range("C5").Value2 = 5.060307348339624
$ws3.Range("D5").Value2 = 4.390548150434426
$ws3.Range("F5").Value2 = -4.644911099880542
$ws3.Range("G5").Value2 = 0.6480446241201373
$ws3.Range("H5").Value2 = 0.9576830352092067
$ws3.Range("I5").Value2 = 1.563625082763474
$ws3.Range("J5").Value2 = -4.708006675992304
$ws3.Range("B6").Value2 = -3.334758169096592
$ws3.Range("C6").Value2 = 6.162846263287839
$ws3.Range("D6").Value2 = 5.950397803041026
$ws3.Range("E6").Value2 = 4.644911099880542
$ws3.Range("G6").Value2 = 5.649320023779967
$ws3.Range("H6").Value2 = 5.879268145723852
$ws3.Range("I6").Value2 = 5.323877283435174
$ws3.Range("J6").Value2 = 3.285586185652335
$ws3.Range("B7").Value2 = -6.386564016778547
$ws3.Range("C7").Value2 = 3.771893612103899
$ws3.Range("D7").Value2 = 3.0074511223745
$ws3.Range("E7").Value2 = -0.6480446241201373
$ws3.Range("F7").Value2 = -5.649320023779967
$ws3.Range("H7").Value2 = 0.855178496583057
$ws3.Range("I7").Value2 = 1.120352506054777
$ws3.Range("J7").Value2 = -8.545146916659615
$ws3.Range("B8").Value2 = -6.530861996158466
$ws3.Range("C8").Value2 = 3.866228331476032
$ws3.Range("D8").Value2 = 3.081890267858277
$ws3.Range("E8").Value2 = -0.9576830352092067
$ws3.Range("F8").Value2 = -5.879268145723852
$ws3.Range("G8").Value2 = -0.855178496583057
$ws3.Range("I8").Value2 = 0.9582870329572257
$ws3.Range("J8").Value2 = -14.21979675563483
$ws3.Range("B9").Value2 = -5.54463229587341
$ws3.Range("C9").Value2 = 1.995625243972477
$ws3.Range("D9").Value2 = 1.098328329486271
$ws3.Range("E9").Value2 = -1.563625082763474
$ws3.Range("F9").Value2 = -5.323877283435174
$ws3.Range("G9").Value2 = -1.120352506054777
$ws3.Range("H9").Value2 = -0.9582870329572257
$ws3.Range("J9").Value2 = -4.017029546964357
$ws3.Range("B10").Value2 = -4.011492991243025
$ws3.Range("C10").Value2 = 8.741976627948864
$ws3.Range("D10").Value2 = 8.472778910051924
$ws3.Range("E10").Value2 = 4.708006675992304
$ws3.Range("F10").Value2 = -3.285586185652335
$ws3.Range("G10").Value2 = 8.545146916659615
$ws3.Range("H10").Value2 = 14.21979675563483
$ws3.Range("I10").Value2 = 4.017029546964357

# --- Resumen ---
$ws4.Range("B2").Value2 = 5
$ws4.Range("D2").Value2 = 3
$ws4.Range("E2").Value2 = 62.5
$ws4.Range("B3").Value2 = 4
$ws4.Range("D3").Value2 = 3
$ws4.Range("E3").Value2 = 50
$ws4.Range("C6").Value2 = 0
$ws4.Range("D6").Value2 = 5
$ws4.Range("C7").Value2 = 0
$ws4.Range("D7").Value2 = 5
$ws4.Range("B8").Value2 = 1
$ws4.Range("D8").Value2 = 1
$ws4.Range("E8").Value2 = 12.5
$ws4.Range("A9").Value2 = "Block Bootstrapping"
$ws4.Range("B9").Value2 = 0
$ws4.Range("D9").Value2 = 1
$ws4.Range("E9").Value2 = 0
$ws4.Range("F9").Value2 = 7.248884056718633
$ws4.Range("A10").Value2 = "AREPD"
$ws4.Range("C10").Value2 = 6
$ws4.Range("D10").Value2 = 2
$ws4.Range("F10").Value2 = 6.489367220625381

